$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: insert 2 new rows before row 9 (shifts old rows 9+ down by 2) ---
$ws.Rows("9:10").Insert()

# --- Write cells in the exact order that introduces brand-new shared strings in
#     the sequence: T003, T004, product1, product2, product3, T005,
#     catalogDao.* (row 5-10), catalogService.* (row 17-21) -- this reproduces the
#     shared-strings table ordering from the target workbook. ---

$ws.Range("C8").Value = "T003"
$ws.Range("C9").Value = "T004"
$ws.Range("E6").Value = "product1"
$ws.Range("E7").Value = "product2"
$ws.Range("E8").Value = "product3"
$ws.Range("C10").Value = "T005"

$ws.Range("D5").Value = "catalogDao.findAll()"
$ws.Range("D6").Value = "catalogDao.findByCode(""P001"")"
$ws.Range("D7").Value = "catalogDao.findByCode(""P002"")"
$ws.Range("D8").Value = "catalogDao.findByCode(""P003"")"
$ws.Range("D9").Value = "catalogDao.findByCode(""P111"")"
$ws.Range("D10").Value = "catalogDao.findByCode(""P112"")"

$ws.Range("D17").Value = "catalogService.listAllProducts()"
$ws.Range("D18").Value = "catalogService.getProductByCode(""P002"")"
$ws.Range("D19").Value = "catalogService.getProductByCode(""P001"")"
$ws.Range("D20").Value = "catalogService.getProductByCode(""P111"")"
$ws.Range("D21").Value = "catalogService.getProductByCode(""P112"")"

# --- Remaining cells: A/B columns + already-existing strings reused in E/C ---

$ws.Range("B5").Value = "ALL_PRODUCT"
$ws.Range("C5").Value = "T001"
$ws.Range("E5").Value = "List<Product>"

$ws.Range("B6").Value = "PRODUCT_BY_CODE"
$ws.Range("C6").Value = "T001"

$ws.Range("B7").Value = "PRODUCT_BY_CODE"
$ws.Range("C7").Value = "T002"

$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "PRODUCT_BY_CODE"

$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "PRODUCT_BY_CODE"
$ws.Range("E9").Value = "no product"

$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "PRODUCT_BY_CODE"
$ws.Range("E10").Value = "no product"

$ws.Range("E17").Value = "List<Product>"

$ws.Range("B18").Value = "PRODUCT_BY_CODE"
$ws.Range("E18").Value = "some product"

$ws.Range("B19").Value = "PRODUCT_BY_CODE"
$ws.Range("C19").Value = "T002"
$ws.Range("E19").Value = "some product"

$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "PRODUCT_BY_CODE"
$ws.Range("C20").Value = "T003"
$ws.Range("E20").Value = "no product"

$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "PRODUCT_BY_CODE"
$ws.Range("C21").Value = "T004"
$ws.Range("E21").Value = "no product"

# --- Column D width + selection ---
$ws.Columns("D").ColumnWidth = 39
$ws.Range("G12").Select()
